$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Widen column A (raw OOXML width of 23 corresponds to a COM ColumnWidth
# of 23 - 5/6 due to the default font padding offset applied by Excel)
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668

# Update car names to include model + model year suffix
$ws.Range("A2").Value = "BYD_SEAL_2023_"
$ws.Range("A3").Value = "XPENG_G9_2023_"
$ws.Range("A4").Value = "VinFast_VF8_2023_"
$ws.Range("A5").Value = "Honda_ZR-V_2023_"
$ws.Range("A6").Value = "BYD_SEAL-U_2023_"
$ws.Range("A7").Value = "Volkswagen_ID.7_2023_"
$ws.Range("A8").Value = "BMW_5 Series_2023_"
$ws.Range("A9").Value = "smart_#3_"
$ws.Range("A10").Value = "BYD_Tang_2023_"
$ws.Range("A11").Value = "Hyundai_KONA_2023_"
$ws.Range("A12").Value = "Kia_EV9_2023_"
$ws.Range("A13").Value = "NIO_ET5_2023_"
$ws.Range("A14").Value = "NIO_EL7_2023_"
$ws.Range("A15").Value = "Lexus_RZ_2023_"
